$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.191.12"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.658.65"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'0.5187"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("D9").Value = "'0.06272"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").Value = "'20.77"
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("D11").Value = "'0.07776"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'4.473"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "1.635.65"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").Value = "1.886.70"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'0.5468"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").Value = "0.0₅8117"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").Value = "'64.96"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "26.197.40"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "'4.615"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").Value = "'191.85"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "'10.07"
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("D23").Value = "'6.004"
$ws.Range("E23").Value = "  -5.02%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'139.43"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'0.1221"
$ws.Range("E26").Value = "  -3.69%  "
$ws.Range("D27").Value = "'7.291"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'16.15"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "'1.437"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").Value = "'0.05935"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").Value = "'3.549"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").Value = "'3.276"
$ws.Range("D34").Value = "'1.581"
$ws.Range("E34").Value = "  -6.09%  "
$ws.Range("D35").Value = "'0.9604"
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("D36").Value = "'2.420"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "'2.768"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "'0.5688"
$ws.Range("E38").Value = "  -6.68%  "
$ws.Range("D39").Value = "'6.035"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "'0.01592"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").Value = "'0.8510"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "1.007.93"
$ws.Range("E43").Value = "  -7.92%  "
$ws.Range("D44").Value = "'100.60"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "1.800.86"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "'56.48"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").Value = "'8.009"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").Value = "'0.4281"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").Value = "'0.05168"
$ws.Range("E51").Value = "  -0.76%  "
